# Updated cryptos list on Sat Mar 25 18:46:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.581.82'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.748.65'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.52'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4585'
$ws.Range("E7").Value = '  +8.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3538'
$ws.Range("E8").Value = '  -2.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07459'
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.04'
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.091'
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.67'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.974'
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.083'
$ws.Range("E15").Value = '  -2.51%  '
$ws.Range("D16").Value = '1.745.41'
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.12'
$ws.Range("E17").Value = '  +1.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001061'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06407'
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.72'
$ws.Range("E21").Value = '  -1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.767'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("D23").Value = '27.636.20'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.119'
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("E26").Value = '  +4.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.10'
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").Value = '1.947.60'
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.058'
$ws.Range("E29").Value = '  -3.12%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.64'
$ws.Range("E30").Value = '  +1.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.052'
$ws.Range("E31").Value = '  -6.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09214'
$ws.Range("E32").Value = '  +4.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.663'
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.514'
$ws.Range("E34").Value = '  -0.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02286'
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("E36").Value = '  -4.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06016'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.945'
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6284'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.198'
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.378'
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.738'
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.09'
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.721'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.49'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.932'
$ws.Range("E48").Value = '  -1.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06857'
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.128'
$ws.Range("E50").Value = '  -3.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.61'
$ws.Range("E51").Value = '  -2.36%  '
